$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-06 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-07 Thursday", 2) | Out-Null
$d.Content.Find.Execute("83×21=1743", $true, $false, $false, $false, $false, $true, 1, $false, "37×55=2035", 2) | Out-Null
$d.Content.Find.Execute("89×94=8366", $true, $false, $false, $false, $false, $true, 1, $false, "82×70=5740", 2) | Out-Null
$d.Content.Find.Execute("61×37=2257", $true, $false, $false, $false, $false, $true, 1, $false, "23×12=276", 2) | Out-Null
$d.Content.Find.Execute("26×67=1742", $true, $false, $false, $false, $false, $true, 1, $false, "82×30=2460", 2) | Out-Null
$d.Content.Find.Execute("64×36=2304", $true, $false, $false, $false, $false, $true, 1, $false, "20×72=1440", 2) | Out-Null
$d.Content.Find.Execute("65×74=4810", $true, $false, $false, $false, $false, $true, 1, $false, "81×56=4536", 2) | Out-Null
$d.Content.Find.Execute("71×77=5467", $true, $false, $false, $false, $false, $true, 1, $false, "84×76=6384", 2) | Out-Null
$d.Content.Find.Execute("56×69=3864", $true, $false, $false, $false, $false, $true, 1, $false, "38×53=2014", 2) | Out-Null
$d.Content.Find.Execute("81×88=7128", $true, $false, $false, $false, $false, $true, 1, $false, "97×16=1552", 2) | Out-Null
$d.Content.Find.Execute("63×36=2268", $true, $false, $false, $false, $false, $true, 1, $false, "66×14=924", 2) | Out-Null
$d.Content.Find.Execute("45×81=3645", $true, $false, $false, $false, $false, $true, 1, $false, "20×75=1500", 2) | Out-Null
$d.Content.Find.Execute("53×46=2438", $true, $false, $false, $false, $false, $true, 1, $false, "77×77=5929", 2) | Out-Null
$d.Content.Find.Execute("94×76=7144", $true, $false, $false, $false, $false, $true, 1, $false, "90×32=2880", 2) | Out-Null
$d.Content.Find.Execute("22×41=902", $true, $false, $false, $false, $false, $true, 1, $false, "24×25=600", 2) | Out-Null
$d.Content.Find.Execute("38×72=2736", $true, $false, $false, $false, $false, $true, 1, $false, "20×79=1580", 2) | Out-Null
$d.Content.Find.Execute("17×24=408", $true, $false, $false, $false, $false, $true, 1, $false, "84×44=3696", 2) | Out-Null
$d.Content.Find.Execute("53×34=1802", $true, $false, $false, $false, $false, $true, 1, $false, "67×19=1273", 2) | Out-Null
$d.Content.Find.Execute("20×78=1560", $true, $false, $false, $false, $false, $true, 1, $false, "79×88=6952", 2) | Out-Null
$d.Content.Find.Execute("76×57=4332", $true, $false, $false, $false, $false, $true, 1, $false, "70×91=6370", 2) | Out-Null
$d.Content.Find.Execute("33×65=2145", $true, $false, $false, $false, $false, $true, 1, $false, "28×79=2212", 2) | Out-Null
$d.Content.Find.Execute("34×84=2856", $true, $false, $false, $false, $false, $true, 1, $false, "72×40=2880", 2) | Out-Null
$d.Content.Find.Execute("57×27=1539", $true, $false, $false, $false, $false, $true, 1, $false, "84×82=6888", 2) | Out-Null
$d.Content.Find.Execute("31×29=899", $true, $false, $false, $false, $false, $true, 1, $false, "87×34=2958", 2) | Out-Null
$d.Content.Find.Execute("11×85=935", $true, $false, $false, $false, $false, $true, 1, $false, "98×93=9114", 2) | Out-Null
$d.Content.Find.Execute("39×76=2964", $true, $false, $false, $false, $false, $true, 1, $false, "23×54=1242", 2) | Out-Null
